$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.931.93'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '2.702.49'
$ws.Range('E3').Value = '  +1.81%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '608.88'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.19%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '158.37'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.18%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.589'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.75%  '
$ws.Range('E9').Value = '  +5.11%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.07'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +4.71%  '
$ws.Range('E12').Value = '  +1.28%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '30.22'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +3.90%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000205'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +10.66%  '
$ws.Range('D15').Value = '3.185.65'
$ws.Range('E15').Value = '  +1.65%  '
$ws.Range('D16').Value = '65.766.20'
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('D17').Value = '2.701.61'
$ws.Range('E17').Value = '  +0.92%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.79'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('E19').Value = '  +1.35%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '360.97'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.48%  '
$ws.Range('E21').Value = '  +3.59%  '
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '70.23'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.85%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.84'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +3.28%  '
$ws.Range('E25').Value = '  +11.85%  '
$ws.Range('B26').Value = 'SuiNetwork'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.64'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.51%  '
$ws.Range('B27').Value = 'Fetch.AI'
$ws.Range('C27').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.70'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +3.42%  '
$ws.Range('E28').Value = '  +4.32%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.32'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.87%  '
$ws.Range('E30').Value = '  +4.86%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.998'
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '535.19'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.85%  '
$ws.Range('E34').Value = '  +5.91%  '
$ws.Range('E35').Value = '  -2.63%  '
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '20.81'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.27%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '162.69'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.37%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.01'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.65%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '42.82'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.53%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '169.06'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.94%  '
$ws.Range('E44').Value = '  +1.92%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0620'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '23.68'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +2.92%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.28'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.70%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0268'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +4.40%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.661'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.71%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '21.04'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +7.68%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0987'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.06%  '
